# Updates the cryptos list (prices and 1h volume changes) to the latest
# snapshot pulled by the GitHub Actions scraper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.931.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "'3.258.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'581.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'184.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "'0.129"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.49%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "'0.407"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("D12").Value = "'3.828.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "'27.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").Value = "'67.993.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "'3.257.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").Value = "'13.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").Value = "'415.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.94%  "
$ws.Range("D21").Value = "'7.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.99%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'71.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "'0.507"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").Value = "'9.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").Value = "'22.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("D31").Value = "'5.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.47%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").Value = "'162.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("B36").Value = "Stacks"
$ws.Range("C36").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D36").Value = "'1.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").Value = "'26.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").Value = "'0.794"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'4.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.52%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'6.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'2.634.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'40.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0674"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("D44").Value = "'2.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.25%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'337.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'24.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0273"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'6.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'0.976"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.21%  "
